# fix export of lists in `antibiotics` and `antivirals`
#
# The `synonyms` column stored R-style vector literals, e.g.
#   c("Abacavir", "Abacavir sulfate", "Ziagen")
# These should instead be plain comma-separated strings, e.g.
#   Abacavir,Abacavir sulfate,Ziagen
#
# Walk every used cell on the active sheet and, wherever the text matches
# the c("...", "...") pattern, rewrite it as a comma separated list of the
# quoted items (in order), dropping the c(...) wrapper and the quotes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rowCount = $used.Rows.Count()
$colCount = $used.Columns.Count()

$listPattern = [regex]'^c\((.*)\)$'
$itemPattern = [regex]'"((?:[^"\\]|\\.)*)"'

for ($r = 1; $r -le $rowCount; $r++) {
    for ($c = 1; $c -le $colCount; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value()
        if ($val -ne $null -and $val -is [string]) {
            $m = $listPattern.Match($val)
            if ($m.Success) {
                $inner = $m.Groups[1].Value
                $items = $itemPattern.Matches($inner) | ForEach-Object { $_.Groups[1].Value }
                $joined = [string]::Join(",", $items)
                $cell.Value = $joined
            }
        }
    }
}
